# EIA Table 2.1.A monthly refresh: add November 2016 data, update
# Year-to-Date and Rolling-12-Months totals, and refresh header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the title / subtitle text (October -> November 2016)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Thousand Tons)"

# ---------------------------------------------------------------------
# 2. Insert a new row for "November" under the Year 2016 section
#    (old row 53 "Year to Date" header and everything below shifts
#    down by one row).
# ---------------------------------------------------------------------
$ws.Rows("53:53").Insert()

# Copy the formatting from the preceding data row (October, row 52)
# so the new row matches the existing style (label cell + number cells).
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 48126
$ws.Range("C53").Value = 35274
$ws.Range("D53").Value = 12624
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 215

# ---------------------------------------------------------------------
# 3. Update the "Year to Date" totals (now rows 55-57)
# ---------------------------------------------------------------------
$ws.Range("A55").Value = 2014
$ws.Range("B55").Value = 785834
$ws.Range("C55").Value = 574259
$ws.Range("D55").Value = 207134
$ws.Range("E55").Value = 185
$ws.Range("F55").Value = 4255

$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 689370
$ws.Range("C56").Value = 501628
$ws.Range("D56").Value = 183886
$ws.Range("E56").Value = 149
$ws.Range("F56").Value = 3707

$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 613093
$ws.Range("C57").Value = 450965
$ws.Range("D57").Value = 158865
$ws.Range("E57").Value = 133
$ws.Range("F57").Value = 3130

# ---------------------------------------------------------------------
# 4. Update "Rolling 12 Months Ending in ..." header text and its
#    totals (now rows 58-60)
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("A59").Value = 2015
$ws.Range("B59").Value = 757169
$ws.Range("C59").Value = 551604
$ws.Range("D59").Value = 201320
$ws.Range("E59").Value = 165
$ws.Range("F59").Value = 4081

$ws.Range("A60").Value = 2016
$ws.Range("B60").Value = 663317
$ws.Range("C60").Value = 488843
$ws.Range("D60").Value = 170906
$ws.Range("E60").Value = 147
$ws.Range("F60").Value = 3421
